$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 5 and row 6 (swap their contents)
$cols = @("A", "B", "D", "E", "F", "G", "H", "J", "Q", "R")

foreach ($col in $cols) {
    $cell5 = $ws.Range("$col" + "5")
    $cell6 = $ws.Range("$col" + "6")
    $v5 = $cell5.Value()
    $v6 = $cell6.Value()
    $cell5.Value = $v6
    $cell6.Value = $v5
}
